$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4113.490509519408
$ws.Range("C3").Value = 4113.490509519408
$ws.Range("C4").Value = 3998.138070531252
$ws.Range("C5").Value = 3935.228258331817
$ws.Range("C6").Value = 3935.228258331817
$ws.Range("C7").Value = 3935.228258331817
$ws.Range("C8").Value = 3935.228258331817
$ws.Range("C9").Value = 3935.228258331817
$ws.Range("C10").Value = 3935.228258331817
$ws.Range("C11").Value = 3935.228258331817
$ws.Range("C12").Value = 3935.228258331817
